# Updates the betting odds values in the "Jogos da Semana" worksheet
# (rows 2-7) to match the latest FlashScore snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.18
$ws.Range("H2").Value = 2.55
$ws.Range("J2").Value = 2.95
$ws.Range("K2").Value = 1.72
$ws.Range("L2").Value = 5.2
$ws.Range("M2").Value = 1.19
$ws.Range("N2").Value = 4.15
$ws.Range("O2").Value = 1.75
$ws.Range("P2").Value = 1.95
$ws.Range("Q2").Value = 3.2
$ws.Range("R2").Value = 1.3
$ws.Range("S2").Value = 5.9
$ws.Range("T2").Value = 1.1
$ws.Range("U2").Value = 1.72
$ws.Range("V2").Value = 2
$ws.Range("W2").Value = 2.5
$ws.Range("X2").Value = 1.47
$ws.Range("AA2").Value = 10.25
$ws.Range("AB2").Value = 22
$ws.Range("AD2").Value = 60
$ws.Range("AE2").Value = 4.15
$ws.Range("AF2").Value = 5.6
$ws.Range("AG2").Value = 23
$ws.Range("AH2").Value = 200
$ws.Range("AI2").Value = 7.4
$ws.Range("AK2").Value = 17
$ws.Range("AM2").Value = 75
$ws.Range("AN2").Value = 100

# Row 3
$ws.Range("G3").Value = 2.07
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 2.6
$ws.Range("L3").Value = 3.4
$ws.Range("S3").Value = 2.47
$ws.Range("T3").Value = 1.42
$ws.Range("W3").Value = 1.57
$ws.Range("X3").Value = 2.1
$ws.Range("Y3").Value = 9.25
$ws.Range("Z3").Value = 11
$ws.Range("AA3").Value = 8.75
$ws.Range("AB3").Value = 19.5
$ws.Range("AC3").Value = 15
$ws.Range("AD3").Value = 23
$ws.Range("AE3").Value = 13.5
$ws.Range("AF3").Value = 7.3
$ws.Range("AG3").Value = 13
$ws.Range("AI3").Value = 11.5
$ws.Range("AJ3").Value = 17
$ws.Range("AK3").Value = 10.75
$ws.Range("AL3").Value = 37
$ws.Range("AM3").Value = 23
$ws.Range("AN3").Value = 28
$ws.Range("AO3").Value = 300

# Row 4
$ws.Range("G4").Value = 1.67
$ws.Range("H4").Value = 3.9
$ws.Range("I4").Value = 4.5
$ws.Range("L4").Value = 4.5
$ws.Range("U4").Value = 1.33
$ws.Range("Z4").Value = 9
$ws.Range("AF4").Value = 7.5
$ws.Range("AN4").Value = 34

# Row 5
$ws.Range("S5").Value = 2.1
$ws.Range("T5").Value = 1.67
$ws.Range("U5").Value = 1.22

# Row 6
$ws.Range("N6").Value = 12
$ws.Range("O6").Value = 1.13
$ws.Range("P6").Value = 5.5
$ws.Range("Z6").Value = 7.5
$ws.Range("AA6").Value = 9.5
$ws.Range("AK6").Value = 23

# Row 7
$ws.Range("G7").Value = 2.15
$ws.Range("H7").Value = 3.35
$ws.Range("I7").Value = 3.15
$ws.Range("J7").Value = 2.72
$ws.Range("L7").Value = 3.7
$ws.Range("T7").Value = 1.27
$ws.Range("V7").Value = 2.72
$ws.Range("X7").Value = 1.87
$ws.Range("Y7").Value = 7.2
$ws.Range("Z7").Value = 10
$ws.Range("AA7").Value = 9
$ws.Range("AB7").Value = 20
$ws.Range("AC7").Value = 18
$ws.Range("AD7").Value = 30
$ws.Range("AF7").Value = 6.4
$ws.Range("AI7").Value = 8.75
$ws.Range("AJ7").Value = 15.5
$ws.Range("AK7").Value = 11.25
$ws.Range("AL7").Value = 40
$ws.Range("AM7").Value = 29
$ws.Range("AN7").Value = 40
$ws.Range("AO7").Value = 700
